# Apply the cryptos worksheet updates described by the commit diff.
# Each row's Coin/Link/Price/Volume(1h) cells are refreshed to their new
# values; a few rows were also re-ranked, which swaps which coin's data
# occupies that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column always holds plain text in this sheet (values such as
# '53.775.33' or '1.00' are not real numbers). Force every Price cell that
# we touch to Text format first so Excel doesn't silently reinterpret the
# new value as a number (which would drop formatting like trailing zeros).
$priceRows = @(2,3,4,5,6,7,8,9,10,12,13,14,15,16,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '53.775.33'
$ws.Range("E2").Value = '  -4.00%  '

$ws.Range("D3").Value = '2.210.51'
$ws.Range("E3").Value = '  -6.34%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '489.52'
$ws.Range("E5").Value = '  -2.41%  '

$ws.Range("D6").Value = '125.59'
$ws.Range("E6").Value = '  -2.64%  '

$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.48%  '

$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  -3.72%  '

$ws.Range("D9").Value = '2.237.29'
$ws.Range("E9").Value = '  -5.33%  '

$ws.Range("D10").Value = '0.0926'
$ws.Range("E10").Value = '  -5.76%  '

$ws.Range("E11").Value = '  -1.19%  '

$ws.Range("D12").Value = '4.69'
$ws.Range("E12").Value = '  -1.89%  '

$ws.Range("D13").Value = '0.316'
$ws.Range("E13").Value = '  -2.43%  '

$ws.Range("D14").Value = '2.602.73'
$ws.Range("E14").Value = '  -6.32%  '

$ws.Range("D15").Value = '21.23'
$ws.Range("E15").Value = '  -0.31%  '

$ws.Range("D16").Value = '53.707.01'
$ws.Range("E16").Value = '  -4.06%  '

$ws.Range("E17").Value = '  -2.59%  '

$ws.Range("D18").Value = '2.215.26'
$ws.Range("E18").Value = '  -6.62%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '3.98'
$ws.Range("E19").Value = '  -0.85%  '

$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '9.61'
$ws.Range("E20").Value = '  -3.68%  '

$ws.Range("D21").Value = '295.85'
$ws.Range("E21").Value = '  -3.57%  '

$ws.Range("D22").Value = '6.19'
$ws.Range("E22").Value = '  -1.26%  '

$ws.Range("D23").Value = '0.996'
$ws.Range("E23").Value = '  -0.33%  '

$ws.Range("D24").Value = '63.11'
$ws.Range("E24").Value = '  -4.58%  '

$ws.Range("D25").Value = '0.997'
$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("D26").Value = '0.368'
$ws.Range("E26").Value = '  +1.03%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.146'
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.307.33'
$ws.Range("E28").Value = '  -6.70%  '

$ws.Range("D29").Value = '7.06'
$ws.Range("E29").Value = '  -2.02%  '

$ws.Range("D30").Value = '165.95'
$ws.Range("E30").Value = '  -4.10%  '

$ws.Range("D31").Value = '1.59'
$ws.Range("E31").Value = '  -3.04%  '

$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0₃0677'
$ws.Range("E32").Value = '  -4.21%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.997'
$ws.Range("E33").Value = '  -0.25%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '5.79'
$ws.Range("E34").Value = '  +0.29%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '0.992'
$ws.Range("E35").Value = '  -0.59%  '

$ws.Range("E36").Value = '  -0.52%  '

$ws.Range("D37").Value = '17.41'
$ws.Range("E37").Value = '  -0.94%  '

$ws.Range("D38").Value = '1.17'
$ws.Range("E38").Value = '  +0.51%  '

$ws.Range("D39").Value = '0.859'
$ws.Range("E39").Value = '  +6.96%  '

$ws.Range("D40").Value = '3.60'
$ws.Range("E40").Value = '  -2.56%  '

$ws.Range("D41").Value = '35.08'
$ws.Range("E41").Value = '  -3.14%  '

$ws.Range("D42").Value = '0.369'
$ws.Range("E42").Value = '  +0.77%  '

$ws.Range("D43").Value = '1.38'
$ws.Range("E43").Value = '  +0.27%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  -1.27%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '4.84'
$ws.Range("E45").Value = '  +3.56%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '125.75'
$ws.Range("E46").Value = '  -1.84%  '

$ws.Range("D47").Value = '0.0882'
$ws.Range("E47").Value = '  -1.10%  '

$ws.Range("D48").Value = '0.539'
$ws.Range("E48").Value = '  -4.13%  '

$ws.Range("D49").Value = '236.88'
$ws.Range("E49").Value = '  -0.48%  '

$ws.Range("D50").Value = '0.0474'
$ws.Range("E50").Value = '  -1.15%  '

$ws.Range("D51").Value = '0.0201'
$ws.Range("E51").Value = '  -2.14%  '

